# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (fund-holdings detail, same layout as
#    the other quarterly sheets) right between "2021-Q4" and "总计".
# 2. Prepend a "2022-Q1" summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $val) {
    # Force a string value even when it "looks like" a number (e.g. "3.50")
    # so trailing zeros / leading zeros survive, matching the source data
    # which stores these columns as text, not numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4)
$newSheet.Name = "2022-Q1"

# Reuse the exact header / row-label formatting from the 2021-Q4 sheet so
# the new sheet's styles match the rest of the workbook.
$q4.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2:A7").Copy()
$newSheet.Range("A2:A7").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$fundRows = @(
    @(0, "006736", "国投瑞银先进制造混合",          "41.36", "92.33", "3.81", "1.5758", 10),
    @(1, "001704", "国投瑞银进宝灵活配置混合",        "33.25", "92.49", "3.50", "1.1638", 10),
    @(2, "161039", "富国中证1000指数增强LOF",        "21.72", "89.03", "0.72", "0.1564", 6),
    @(3, "515760", "华夏中证浙江国资创新发展ETF",      "2.26",  "98.73", "3.93", "0.0888", 7),
    @(4, "770001", "德邦优化灵活配置混合",            "2.49",  "86.80", "3.47", "0.0864", 2),
    @(5, "512190", "浙商汇金中证浙江凤凰行动50ETF",    "0.51",  "98.94", "3.96", "0.0202", 8)
)

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    Set-TextValue $newSheet.Cells.Item($r, 2) $row[1]
    Set-TextValue $newSheet.Cells.Item($r, 3) $row[2]
    Set-TextValue $newSheet.Cells.Item($r, 4) $row[3]
    Set-TextValue $newSheet.Cells.Item($r, 5) $row[4]
    Set-TextValue $newSheet.Cells.Item($r, 6) $row[5]
    Set-TextValue $newSheet.Cells.Item($r, 7) $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) "总计" sheet: add a new leading row for 2022-Q1, pushing the rest
#    down by one.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Grow the formatted block by one row, copying the formatting of the
# last existing data row down onto the new row 6.
$total.Range("A5:D5").Copy()
$total.Range("A6:D6").PasteSpecial(-4122)

$totalRows = @(
    @(0, "2022-Q1", 6, 3.09),
    @(1, "2021-Q4", 8, 2.22),
    @(2, "2021-Q3", 26, 12.9),
    @(3, "2021-Q2", 11, 1.53),
    @(4, "2021-Q1", 7, 0.17)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

Write-Output "2022-Q1 sheet added; total sheet updated"
